$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) are stored as plain text in the
# workbook (many values look like numbers, e.g. "1.00", "0.427", or have
# multiple dots as thousands separators, e.g. "56.779.21"). Force the
# number format to Text before writing so Excel does not silently
# reinterpret these strings as numeric/date values, then restore the
# default "Normal" style so no stray cell formatting is introduced.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("D2").Value = '56.779.21'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '2.989.83'
$ws.Range("E3").Value = '  -3.35%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '498.79'
$ws.Range("E5").Value = '  -3.52%  '
$ws.Range("D6").Value = '134.55'
$ws.Range("E6").Value = '  +3.34%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '0.427'
$ws.Range("E8").Value = '  -2.51%  '
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = '0.350'
$ws.Range("E11").Value = '  -4.38%  '
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '3.499.34'
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("D15").Value = '56.685.70'
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").Value = '2.986.68'
$ws.Range("E17").Value = '  -3.70%  '
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").Value = '12.35'
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").Value = '7.78'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '326.59'
$ws.Range("E21").Value = '  -2.56%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -6.50%  '
$ws.Range("D24").Value = '62.03'
$ws.Range("E24").Value = '  -6.10%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("D27").Value = '0.0₃0890'
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("D29").Value = '6.40'
$ws.Range("E29").Value = '  -3.95%  '
$ws.Range("D30").Value = '6.83'
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("E31").Value = '  -4.80%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '20.32'
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = '1.16'
$ws.Range("E33").Value = '  -6.20%  '
$ws.Range("D34").Value = '155.18'
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("E35").Value = '  -5.56%  '
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("D37").Value = '5.57'
$ws.Range("E37").Value = '  -7.65%  '
$ws.Range("D38").Value = '0.0674'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("B39").Value = 'RenzoRestakedETH'
$ws.Range("C39").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D39").Value = '3.022.28'
$ws.Range("E39").Value = '  -3.40%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").Value = '22.91'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").Value = '36.37'
$ws.Range("E41").Value = '  -9.13%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("E43").Value = '  -5.98%  '
$ws.Range("D44").Value = '2.234.62'
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("E45").Value = '  -5.07%  '
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("E47").Value = '  -7.67%  '
$ws.Range("E48").Value = '  +11.33%  '
$ws.Range("E49").Value = '  +2.82%  '
$ws.Range("D50").Value = '5.72'
$ws.Range("E50").Value = '  -5.39%  '
$ws.Range("D51").Value = '18.97'
$ws.Range("E51").Value = '  -5.34%  '
$ws.Range("D2:E51").Style = "Normal"

